$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.505.16'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.587.56'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -7.52%  '
$ws.Range('D5').Value = '''553.39'
$ws.Range('E5').Value = '  -2.79%  '
$ws.Range('D6').Value = '''140.93'
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '2.599.84'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = '''6.73'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '''0.160'
$ws.Range('E12').Value = '  +5.31%  '
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').Value = '3.042.94'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '59.505.95'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '''23.18'
$ws.Range('E16').Value = '  +5.78%  '
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '2.597.20'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '''4.55'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = '''340.17'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '''10.36'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').Value = '''6.51'
$ws.Range('E22').Value = '  +4.06%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = '''0.479'
$ws.Range('E24').Value = '  +7.34%  '
$ws.Range('D25').Value = '''63.01'
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('D28').Value = '''7.47'
$ws.Range('E28').Value = '  +2.63%  '
$ws.Range('D29').Value = '0.0₃0772'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('D30').Value = '''0.998'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').Value = '''6.15'
$ws.Range('E32').Value = '  +1.85%  '
$ws.Range('D33').Value = '''157.84'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('D35').Value = '''4.12'
$ws.Range('E36').Value = '  +3.27%  '
$ws.Range('D37').Value = '''0.903'
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').Value = '''37.63'
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.48'
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').Value = '''0.842'
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('E41').Value = '  +1.35%  '
$ws.Range('D42').Value = '''288.41'
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('D43').Value = '''136.30'
$ws.Range('E43').Value = '  +8.65%  '
$ws.Range('D44').Value = '''0.996'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '''0.0973'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('D47').Value = '''10.64'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').Value = '''0.0532'
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.983.73'
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''18.64'
$ws.Range('E51').Value = '  +1.17%  '
